$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.262.76"
$ws.Range("E2").Value = "  +1.66%  "
$ws.Range("D3").Value = "2.183.55"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "255.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.12%  "
$ws.Range("E6").Value = "  +1.70%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "68.04"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.75%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.579"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +8.77%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.82"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "58.64"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.23%  "
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.12"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +9.62%  "
$ws.Range("E14").Value = "  +0.53%  "
$ws.Range("D15").Value = "2.508.54"
$ws.Range("E15").Value = "  +0.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.871"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.42%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.50"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.37%  "
$ws.Range("D18").Value = "2.186.11"
$ws.Range("E18").Value = "  +0.84%  "
$ws.Range("D19").Value = "41.205.73"
$ws.Range("E19").Value = "  +1.76%  "
$ws.Range("E20").Value = "  +2.31%  "
$ws.Range("E21").Value = "  +2.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.91"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "232.47"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.67%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.95"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +10.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +23.59%  "
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.54"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "168.84"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.63"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.93%  "
$ws.Range("E32").Value = "  +1.67%  "
$ws.Range("E33").Value = "  +6.84%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.123"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.97%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.45"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "27.37"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +18.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.19"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +11.64%  "
$ws.Range("E38").Value = "  +1.79%  "
$ws.Range("E39").Value = "  +14.27%  "
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.20"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.63%  "
$ws.Range("B41").Value = "Celestia"
$ws.Range("C41").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.53"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +23.50%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.69"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "64.57"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.69%  "
$ws.Range("E44").Value = "  +6.60%  "
$ws.Range("E45").Value = "  +6.30%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.64"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.91%  "
$ws.Range("E47").Value = "  +4.15%  "
$ws.Range("E48").Value = "  +0.50%  "
$ws.Range("E49").Value = "  +5.67%  "
$ws.Range("E50").Value = "  +2.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.28"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.04%  "
